# chore: update Sheets via scheduled runner
# Refreshes market-price-derived columns (H:N) on several Leve-profit
# worksheets with newly pulled values. Only numeric cell values change;
# no formulas, styles, or sheet structure are touched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7921.185
$ws.Range("I19").Value = 470
$ws.Range("J19").Value = 13043.875
$ws.Range("K19").Value = 470
$ws.Range("L19").Value = 13043.875
$ws.Range("M19").Value = -295
$ws.Range("N19").Value = -13393.875
$ws.Range("H28").Value = 12603.333
$ws.Range("I28").Value = 2610.3333
$ws.Range("J28").Value = 27592.834
$ws.Range("K28").Value = 2610.3333
$ws.Range("L28").Value = 27592.834
$ws.Range("M28").Value = -2125.3333
$ws.Range("N28").Value = -28562.834
$ws.Range("H107").Value = 1096.6154
$ws.Range("I107").Value = 1307.2354
$ws.Range("J107").Value = 698.7778
$ws.Range("K107").Value = 1307.2354
$ws.Range("L107").Value = 698.7778
$ws.Range("M107").Value = 612.7646
$ws.Range("N107").Value = -4538.7778
$ws.Range("H116").Value = 2999.1333
$ws.Range("I116").Value = 2958.1
$ws.Range("J116").Value = 3081.2
$ws.Range("K116").Value = 2958.1
$ws.Range("L116").Value = 3081.2
$ws.Range("M116").Value = 483.9000000000001
$ws.Range("N116").Value = -9965.200000000001
$ws.Range("H129").Value = 1406.6552
$ws.Range("J129").Value = 2022.6666
$ws.Range("L129").Value = 6067.9998
$ws.Range("N129").Value = -16067.9998
$ws.Range("H137").Value = 1397.7894
$ws.Range("I137").Value = 1124.6842
$ws.Range("J137").Value = 1670.8948
$ws.Range("K137").Value = 3374.0526
$ws.Range("L137").Value = 5012.6844
$ws.Range("M137").Value = -824.0526
$ws.Range("N137").Value = -10112.6844
$ws.Range("H141").Value = 4267.35
$ws.Range("I141").Value = 2168.4285
$ws.Range("K141").Value = 6505.2855
$ws.Range("M141").Value = -1325.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3661.54
$ws.Range("I61").Value = 2800.8262
$ws.Range("K61").Value = 2800.8262
$ws.Range("M61").Value = -2588.8262
$ws.Range("H136").Value = 3661.54
$ws.Range("I136").Value = 2800.8262
$ws.Range("K136").Value = 8402.4786
$ws.Range("M136").Value = -5852.4786

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 276500
$ws.Range("I107").Value = 367666.66
$ws.Range("K107").Value = 367666.66
$ws.Range("M107").Value = -365746.66
$ws.Range("H134").Value = 1849.2545
$ws.Range("I134").Value = 1534.1163
$ws.Range("J134").Value = 2978.5
$ws.Range("K134").Value = 4602.3489
$ws.Range("L134").Value = 8935.5
$ws.Range("M134").Value = -2067.3489
$ws.Range("N134").Value = -14005.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 55001
$ws.Range("J4").Value = 55001
$ws.Range("L4").Value = 55001
$ws.Range("N4").Value = -55225
$ws.Range("H31").Value = 4388.4463
$ws.Range("I31").Value = 1169.1613
$ws.Range("J31").Value = 7323.6763
$ws.Range("K31").Value = 1169.1613
$ws.Range("L31").Value = 7323.6763
$ws.Range("M31").Value = -874.1613
$ws.Range("N31").Value = -7913.6763
$ws.Range("H34").Value = 4388.4463
$ws.Range("I34").Value = 1169.1613
$ws.Range("J34").Value = 7323.6763
$ws.Range("K34").Value = 1169.1613
$ws.Range("L34").Value = 7323.6763
$ws.Range("M34").Value = -967.1613
$ws.Range("N34").Value = -7727.6763
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 108.36842
$ws.Range("I2").Value = 84.2
$ws.Range("J2").Value = 117
$ws.Range("K2").Value = 505.2
$ws.Range("L2").Value = 702
$ws.Range("M2").Value = -392.2
$ws.Range("N2").Value = -928
$ws.Range("H58").Value = 1760.8695
$ws.Range("I58").Value = 625
$ws.Range("K58").Value = 1875
$ws.Range("M58").Value = -1747
$ws.Range("H75").Value = 1950
$ws.Range("I75").Value = 900
$ws.Range("J75").Value = 3000
$ws.Range("K75").Value = 2700
$ws.Range("L75").Value = 9000
$ws.Range("M75").Value = -1702
$ws.Range("N75").Value = -10996
$ws.Range("H78").Value = 1950
$ws.Range("I78").Value = 900
$ws.Range("J78").Value = 3000
$ws.Range("K78").Value = 8100
$ws.Range("L78").Value = 27000
$ws.Range("M78").Value = -3108
$ws.Range("N78").Value = -36984
$ws.Range("H81").Value = 6322.077
$ws.Range("I81").Value = 1532.6
$ws.Range("J81").Value = 9315.5
$ws.Range("K81").Value = 4597.799999999999
$ws.Range("L81").Value = 27946.5
$ws.Range("M81").Value = -3474.799999999999
$ws.Range("N81").Value = -30192.5
$ws.Range("H84").Value = 6322.077
$ws.Range("I84").Value = 1532.6
$ws.Range("J84").Value = 9315.5
$ws.Range("K84").Value = 13793.4
$ws.Range("L84").Value = 83839.5
$ws.Range("M84").Value = -8177.4
$ws.Range("N84").Value = -95071.5
$ws.Range("H131").Value = 4933.3335
$ws.Range("J131").Value = 5907.407
$ws.Range("L131").Value = 17722.221
$ws.Range("N131").Value = -27802.221
$ws.Range("H133").Value = 12086.808
$ws.Range("I133").Value = 7031.1113
$ws.Range("J133").Value = 14763.353
$ws.Range("K133").Value = 21093.3339
$ws.Range("L133").Value = 44290.05899999999
$ws.Range("M133").Value = -16033.3339
$ws.Range("N133").Value = -54410.05899999999
$ws.Range("H136").Value = 4412.7144
$ws.Range("I136").Value = 1617.8
$ws.Range("J136").Value = 11400
$ws.Range("K136").Value = 4853.4
$ws.Range("L136").Value = 34200
$ws.Range("M136").Value = 246.6000000000004
$ws.Range("N136").Value = -44400
$ws.Range("H137").Value = 31583.578
$ws.Range("J137").Value = 59106.277
$ws.Range("L137").Value = 177318.831
$ws.Range("N137").Value = -187518.831
$ws.Range("H139").Value = 326224.2
$ws.Range("I139").Value = 626692.5600000001
$ws.Range("J139").Value = 5724.6
$ws.Range("K139").Value = 1880077.68
$ws.Range("L139").Value = 17173.8
$ws.Range("M139").Value = -1874937.68
$ws.Range("N139").Value = -27453.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 8321.666999999999
$ws.Range("J123").Value = 8321.666999999999
$ws.Range("L123").Value = 8321.666999999999
$ws.Range("N123").Value = -13221.667
$ws.Range("H132").Value = 3634.3704
$ws.Range("I132").Value = 2855.2
$ws.Range("K132").Value = 8565.599999999999
$ws.Range("M132").Value = -6035.599999999999
$ws.Range("H136").Value = 7564.4
$ws.Range("J136").Value = 7564.4
$ws.Range("L136").Value = 22693.2
$ws.Range("N136").Value = -27793.2
$ws.Range("H141").Value = 48936.332
$ws.Range("J141").Value = 48936.332
$ws.Range("L141").Value = 48936.332
$ws.Range("N141").Value = -59296.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2488.647
$ws.Range("I82").Value = 2733.8
$ws.Range("J82").Value = 2138.4285
$ws.Range("K82").Value = 2733.8
$ws.Range("L82").Value = 2138.4285
$ws.Range("M82").Value = -2372.8
$ws.Range("N82").Value = -2860.4285
$ws.Range("H85").Value = 2488.647
$ws.Range("I85").Value = 2733.8
$ws.Range("J85").Value = 2138.4285
$ws.Range("K85").Value = 2733.8
$ws.Range("L85").Value = 2138.4285
$ws.Range("M85").Value = -1485.8
$ws.Range("N85").Value = -4634.4285
$ws.Range("H132").Value = 3237.84
$ws.Range("I132").Value = 2610.611
$ws.Range("J132").Value = 4850.7144
$ws.Range("K132").Value = 7831.833
$ws.Range("L132").Value = 14552.1432
$ws.Range("M132").Value = -5301.833
$ws.Range("N132").Value = -19612.1432
$ws.Range("H136").Value = 2668.6667
$ws.Range("I136").Value = 2003
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 6009
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -3459
$ws.Range("N136").Value = -17100
